# The workbook's rows got re-ordered (species/observation records were
# reshuffled while the shared "location batch" columns stayed put).
# For each group of rows below, the values in the "observation-specific"
# columns (Id, Taxonsorteringsordning, Rödlistade, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Aktivitet, Ost, Nord) are reversed in order
# across the rows of the group - i.e. row[0] swaps with row[last], etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (1-based) that carry the per-observation data which moved.
$cols = @(1, 2, 4, 5, 6, 7, 8, 13, 17, 18)

# Row groups whose observation data got reversed in order.
$groups = @(
    @(2, 3, 4, 5),
    @(16, 17),
    @(20, 21),
    @(23, 24)
)

foreach ($rows in $groups) {
    $n = $rows.Count

    # Snapshot current values (Value2 avoids date/currency reinterpretation)
    # for every row/col in this group before writing anything back.
    $snapshot = @{}
    foreach ($r in $rows) {
        foreach ($c in $cols) {
            $snapshot["$r`_$c"] = $ws.Cells.Item($r, $c).Value2
        }
    }

    # Write the snapshot back in reverse row order.
    for ($i = 0; $i -lt $n; $i++) {
        $srcRow = $rows[$i]
        $dstRow = $rows[$n - 1 - $i]
        foreach ($c in $cols) {
            $ws.Cells.Item($dstRow, $c).Value = $snapshot["$srcRow`_$c"]
        }
    }
}
